$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.034748368925986)
    3 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    5 = @(1.505614041169197, 0.3375848360084654, 0.1529057820181812, 6.48142807727062, 0, 8.477532736466465)
    6 = @(0.06328177979961902, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 2.369310727790667)
    7 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    8 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

$wb.Save()
